$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.516.72'
$ws.Range("E2").Value = '  -1.62%  '
$ws.Range("D3").Value = '2.588.16'
$ws.Range("E3").Value = '  -2.22%  '
$ws.Range("E4").Value = '  -0.31%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '561.91'
$ws.Range("E5").Value = '  -1.04%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.02'
$ws.Range("E6").Value = '  -2.75%  '
$ws.Range("E7").Value = '  +0.14%  '
$ws.Range("E8").Value = '  -1.64%  '
$ws.Range("D9").Value = '2.597.42'
$ws.Range("E9").Value = '  -2.82%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.67'
$ws.Range("E10").Value = '  -2.89%  '
$ws.Range("E11").Value = '  -0.72%  '
$ws.Range("E12").Value = '  +10.53%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.356'
$ws.Range("E13").Value = '  +3.90%  '
$ws.Range("D14").Value = '3.041.95'
$ws.Range("E14").Value = '  -2.61%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '23.32'
$ws.Range("E15").Value = '  +5.92%  '
$ws.Range("D16").Value = '59.439.27'
$ws.Range("E16").Value = '  -1.79%  '
$ws.Range("E17").Value = '  -0.21%  '
$ws.Range("D18").Value = '2.599.61'
$ws.Range("E18").Value = '  -2.75%  '
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '338.99'
$ws.Range("E20").Value = '  -1.49%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.41'
$ws.Range("E21").Value = '  -0.72%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.52'
$ws.Range("E22").Value = '  +2.17%  '
$ws.Range("E23").Value = '  +0.23%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.40'
$ws.Range("E24").Value = '  -4.97%  '
$ws.Range("E25").Value = '  +7.52%  '
$ws.Range("E26").Value = '  +0.38%  '
$ws.Range("E27").Value = '  -1.98%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.47'
$ws.Range("E28").Value = '  +0.44%  '
$ws.Range("D29").Value = '0.0₃0780'
$ws.Range("E29").Value = '  -3.04%  '
$ws.Range("E30").Value = '  +0.04%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.23'
$ws.Range("E31").Value = '  -0.84%  '
$ws.Range("E32").Value = '  -2.26%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '158.46'
$ws.Range("E33").Value = '  +1.63%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '19.08'
$ws.Range("E34").Value = '  -0.99%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.06'
$ws.Range("E35").Value = '  -1.16%  '
$ws.Range("E36").Value = '  +0.90%  '
$ws.Range("E37").Value = '  -1.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.871'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '37.45'
$ws.Range("E40").Value = '  -2.06%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '295.40'
$ws.Range("E41").Value = '  -2.76%  '
$ws.Range("E42").Value = '  +0.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '139.75'
$ws.Range("E43").Value = '  +8.34%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.998'
$ws.Range("E44").Value = '  +0.44%  '
$ws.Range("E45").Value = '  -0.92%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.596'
$ws.Range("E46").Value = '  -1.91%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.65'
$ws.Range("E47").Value = '  -0.29%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0532'
$ws.Range("E48").Value = '  -3.07%  '
$ws.Range("E49").Value = '  -0.56%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.79'
$ws.Range("E50").Value = '  -1.08%  '
$ws.Range("D51").Value = '1.967.67'
$ws.Range("E51").Value = '  -0.08%  '
